$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. "1.00", "0.578") that must
# stay stored as TEXT (matches the source inlineStr cells, preserving trailing zeros).
# Force text format first, assign the value, then restore the default "Normal" style
# so no stray quote-prefix/number-format style sticks to the cell.
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D15", "D16", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "524.94"
$ws.Range("D6").Value = "153.72"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "0.578"
$ws.Range("D10").Value = "0.108"
$ws.Range("D11").Value = "0.348"
$ws.Range("D15").Value = "21.82"
$ws.Range("D16").Value = "0.0000141"
$ws.Range("D18").Value = "4.73"
$ws.Range("D19").Value = "352.25"
$ws.Range("D20").Value = "10.59"
$ws.Range("D21").Value = "6.25"
$ws.Range("D23").Value = "61.23"
$ws.Range("D24").Value = "0.427"
$ws.Range("D25").Value = "0.166"
$ws.Range("D26").Value = "0.997"
$ws.Range("D28").Value = "7.23"
$ws.Range("D29").Value = "1.00"
$ws.Range("D30").Value = "6.14"
$ws.Range("D31").Value = "1.61"
$ws.Range("D32").Value = "19.30"
$ws.Range("D33").Value = "149.89"
$ws.Range("D34").Value = "4.08"
$ws.Range("D35").Value = "1.18"
$ws.Range("D36").Value = "0.902"
$ws.Range("D37").Value = "0.886"
$ws.Range("D38").Value = "36.82"
$ws.Range("D39").Value = "304.47"
$ws.Range("D40").Value = "3.74"
$ws.Range("D41").Value = "1.46"
$ws.Range("D42").Value = "0.636"
$ws.Range("D44").Value = "20.08"
$ws.Range("D45").Value = "0.997"
$ws.Range("D46").Value = "0.0557"
$ws.Range("D47").Value = "0.0239"
$ws.Range("D48").Value = "4.82"
$ws.Range("D49").Value = "10.35"
$ws.Range("D50").Value = "18.98"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining cells (coin names, links, already-text prices, percentage strings)
# can be assigned directly; Excel keeps them as text since they are not valid numbers.
$ws.Range("D2").Value = "60.614.03"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.643.86"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  -3.48%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "3.108.85"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "60.600.29"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.654.39"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "0.0₃0845"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -2.12%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E32").Value = "  -1.04%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +3.84%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").Value = "1.980.79"
$ws.Range("E51").Value = "  -0.82%  "
